$wb = $excel.ActiveWorkbook

# --- SQL sheet: append a new row (13) ---
$sql = $wb.Worksheets.Item("SQL")
$sql.Cells.Item(13, 1).Value = 1162
$sql.Cells.Item(13, 2).Value = "Тестовое сообщение"
$sql.Cells.Item(13, 3).Value = "Тестовое сообщение"
$sql.Cells.Item(13, 4).Value = "Тестовое сообщение"

# --- Python sheet: append a new row (30) ---
$py = $wb.Worksheets.Item("Python")
$py.Cells.Item(30, 1).Value = 2090
$py.Cells.Item(30, 2).Value = "Тестовое сообщение"
$py.Cells.Item(30, 3).Value = "Тестовое сообщение"
$py.Cells.Item(30, 4).Value = "Тестовое сообщение1"
